$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11338.2597232856
$ws.Range("C2").Value = 10473.8973076142
$ws.Range("E2").Value = 6959.6259946357
$ws.Range("F2").Value = 258.385970927079

$ws.Range("B3").Value = 11375.6392876749
$ws.Range("C3").Value = 10492.7196954354
$ws.Range("E3").Value = 6834.95583953076
$ws.Range("F3").Value = 253.975647290257

$ws.Range("B4").Value = 11444.8654743131
$ws.Range("C4").Value = 10538.3436487851
$ws.Range("E4").Value = 7138.92239066545
$ws.Range("F4").Value = 268.541918310439

$ws.Range("B5").Value = 11391.824268248
$ws.Range("C5").Value = 9847.00963593433
$ws.Range("E5").Value = 7111.71769462548
$ws.Range("F5").Value = 238.602805439992

$ws.Range("B6").Value = 4296.32131239475
$ws.Range("C6").Value = 6896.23776963839
$ws.Range("E6").Value = 6699.90339760516
$ws.Range("F6").Value = 98.4950486351481

$ws.Range("B7").Value = 3918.1333122187
$ws.Range("C7").Value = 6738.67220508619
$ws.Range("E7").Value = 6425.93317393359
$ws.Range("F7").Value = 80.5143907924907
